$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 216
$ws.Range("J2").Value = 778
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 235
$ws.Range("M2").Value = 18
$ws.Range("N2").Value = 152
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = 10
$ws.Range("S2").Value = 83
$ws.Range("T2").Value = 134
$ws.Range("U2").Value = 10
$ws.Range("V2").Value = 1268
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1172
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 17
$ws.Range("AA2").Value = 9
